# edit.ps1 - apply the darkness.docx power-table edits described in the
# commit "Dropped armor values a little. Making list of traits / maneuvers"
#
# Strategy: Word's Range.InsertXML replaces the *entire paragraph* that
# contains the Find-matched range (it is not a surgical sub-range patch),
# so each edit below locates its target paragraph with Find, then feeds
# InsertXML the complete original paragraph XML with only the intended
# runs/bookmarks changed - this preserves every untouched attribute
# (w14:paraId, rsids, pPr, sibling runs, etc.) exactly as Word would.

$d = $word.ActiveDocument

function Apply-ParagraphXml($FindText, $ParagraphXml) {
    $rng = $d.Content
    $found = $rng.Find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find target text: $FindText"
    }

    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' + $ParagraphXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $rng.InsertXML($pkg)
}

# ---------------------------------------------------------------------
# 1) "E1" -> "M"  (Burrow/Crawl-ish speed cell next to "Self")
# ---------------------------------------------------------------------
$p1 = '<w:p w14:paraId="490496BB" w14:textId="5844C798" w:rsidR="00A4187C" w:rsidRPr="00F53F89" w:rsidRDefault="001D1EB9"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr></w:pPr><w:r w:rsidRPr="00F53F89"><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t>M</w:t></w:r></w:p>'
Apply-ParagraphXml "E1" $p1

# ---------------------------------------------------------------------
# 2) "Armor 0/6/6" -> "Armor 0/" + "4" + "/" + "4" (separate runs)
# ---------------------------------------------------------------------
$p2 = '<w:p w14:paraId="1F89CE77" w14:textId="26F8D08B" w:rsidR="00A4187C" w:rsidRPr="00F53F89" w:rsidRDefault="00967876" w:rsidP="00885C8B"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:line="276" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr></w:pPr><w:r w:rsidRPr="00F53F89"><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t>Armor 0/</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t>4</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t>/</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t>4</w:t></w:r></w:p>'
Apply-ParagraphXml "Armor 0/6/6" $p2

# ---------------------------------------------------------------------
# 3) "Armor 0/2/2 / x3 / -- / " -> "Armor" + _GoBack bookmark + " / x3 / -- / "
#    (the trailing "10P" / " ***" runs in the same paragraph are kept as-is)
# ---------------------------------------------------------------------
$p3 = '<w:p w14:paraId="04C2647B" w14:textId="42D3BCD9" w:rsidR="00A4187C" w:rsidRPr="00F53F89" w:rsidRDefault="00845D13" w:rsidP="001A239F"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr></w:pPr><w:r w:rsidRPr="00F53F89"><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t>Armor</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t xml:space="preserve"> / x3 / -- / </w:t></w:r><w:r w:rsidR="006D0CDD" w:rsidRPr="0095224D"><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t>10P</w:t></w:r><w:r w:rsidR="006D0CDD"><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t xml:space="preserve"> ***</w:t></w:r></w:p>'
Apply-ParagraphXml "Armor 0/2/2 / x3 / -- / " $p3

# ---------------------------------------------------------------------
# 4) Drop the old _GoBack bookmark that used to sit after "Move / x3 / +1B / 10P"
#    (it moved to the Armor paragraph above, so this one is now just removed)
# ---------------------------------------------------------------------
$p4 = '<w:p w14:paraId="1910E555" w14:textId="6639871A" w:rsidR="006F0A18" w:rsidRPr="00F53F89" w:rsidRDefault="006F0A18" w:rsidP="00845D13"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t>Move / x3 / +1B / 10P</w:t></w:r></w:p>'
Apply-ParagraphXml "Move / x3 / +1B / 10P" $p4

Write-Output "Done applying edits."
